# Horarios actualizados Línea 141 - 371
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912": new scrape arrived (06:57:11), one new stop was
# inserted right after the current last-but-one batch (row 56),
# and two more stops were appended at the very end (rows 62-63).
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header metadata
$ws1.Range("A2").Value = "Última actualización: 06:57:11"
$ws1.Range("A3").Value = "Total filas: 58"

# Insert the new row at position 56, shifting the previous 56-60 down to 57-61
$ws1.Rows.Item(56).Insert()
$ws1.Range("A56").Value = "06:57:11"
$ws1.Range("B56").Value = "08:05"
$ws1.Range("C56").Value = "23_HERNANDEZ"
$ws1.Range("D56").Value = 68
$ws1.Range("E56").Value = "LP1912"

# Append two brand-new rows at the bottom of the sheet
$ws1.Range("A62").Value = "06:57:11"
$ws1.Range("B62").Value = "08:42"
$ws1.Range("C62").Value = "81_EL PELIGRO"
$ws1.Range("D62").Value = 105
$ws1.Range("E62").Value = "LP1912"

$ws1.Range("A63").Value = "06:57:11"
$ws1.Range("B63").Value = "08:54"
$ws1.Range("C63").Value = "17_ROMERO"
$ws1.Range("D63").Value = 117
$ws1.Range("E63").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "LP1912-215": only the "last updated" timestamp changed.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:57:11"

# ---------------------------------------------------------------
# Sheet "6203-6173": timestamp refreshed and two new stops
# appended at the end (rows 15-16).
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:57:11"
$ws3.Range("A3").Value = "Total filas: 11"

$ws3.Range("A15").Value = "06:57:11"
$ws3.Range("B15").Value = "08:10"
$ws3.Range("C15").Value = "215C_LA PLATA"
$ws3.Range("D15").Value = 73
$ws3.Range("E15").Value = "L6203"

$ws3.Range("A16").Value = "06:57:11"
$ws3.Range("B16").Value = "08:38"
$ws3.Range("C16").Value = "215A_LA PLATA"
$ws3.Range("D16").Value = 101
$ws3.Range("E16").Value = "L6173"

Write-Output "done"
